$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): row -> new F value
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 3091
$wsExhibit.Range("F4").Value = 67
$wsExhibit.Range("F5").Value = 57
$wsExhibit.Range("F6").Value = 13
$wsExhibit.Range("F9").Value = 1076
$wsExhibit.Range("F10").Value = 15092
$wsExhibit.Range("F11").Value = 197
$wsExhibit.Range("F12").Value = 149
$wsExhibit.Range("F13").Value = 509
$wsExhibit.Range("F14").Value = 5986
$wsExhibit.Range("F15").Value = 610
$wsExhibit.Range("F16").Value = 90
$wsExhibit.Range("F18").Value = 95
$wsExhibit.Range("F19").Value = 1250
$wsExhibit.Range("F21").Value = 101
$wsExhibit.Range("F25").Value = 2965
$wsExhibit.Range("F27").Value = 10820
$wsExhibit.Range("F28").Value = 1220
$wsExhibit.Range("F29").Value = 93

# Sheet "全部类型" (sheet4): row -> new F value (same data, offset by +1 row vs 展览)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 3091
$wsAll.Range("F5").Value = 67
$wsAll.Range("F6").Value = 57
$wsAll.Range("F7").Value = 13
$wsAll.Range("F10").Value = 1076
$wsAll.Range("F11").Value = 15092
$wsAll.Range("F12").Value = 197
$wsAll.Range("F13").Value = 149
$wsAll.Range("F14").Value = 509
$wsAll.Range("F15").Value = 5986
$wsAll.Range("F16").Value = 610
$wsAll.Range("F17").Value = 90
$wsAll.Range("F19").Value = 95
$wsAll.Range("F20").Value = 1250
$wsAll.Range("F22").Value = 101
$wsAll.Range("F26").Value = 2965
$wsAll.Range("F29").Value = 10820
$wsAll.Range("F30").Value = 1220
$wsAll.Range("F31").Value = 93
